# Added support for inverse perpetual and added code to validate data in db.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# ListOfValues sheet: add new Pairs values (BTCUSD / ETHUSD) supporting the
# new inverse-perpetual pairs, and leave the selection parked on E6.
# ---------------------------------------------------------------------------
$wsList = $wb.Worksheets.Item("ListOfValues")
$wsList.Range("E4").Value = "BTCUSD"
$wsList.Range("E5").Value = "ETHUSD"
$wsList.Range("E6").Select()

# ---------------------------------------------------------------------------
# Sheet1: switch the sample test row to use the new inverse-perpetual pair
# and a plain MACD strategy, add two validation rows under the header/sample
# row, and trim the now-unused blank row 20 from the bottom of the table.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")

$ws1.Range("C2").Value = "BTCUSD"
$ws1.Range("I2").Value = "MACD"

$ws1.Range("K2").Copy()
$ws1.Range("K3").PasteSpecial(-4122)
$ws1.Range("K4").PasteSpecial(-4122)
$ws1.Rows.Item(3).RowHeight = 14.25
$ws1.Rows.Item(4).RowHeight = 14.25

$ws1.Rows.Item(20).Clear()

# Reset the view: scroll back to column A (removing the old topLeftCell="B1")
# and select the newly inserted row 3 across the table width, restoring
# Sheet1 as the active tab.
$win = $excel.ActiveWindow
$win.ScrollColumn = 1
$win.ScrollRow = 1
$ws1.Range("A3:XFD3").Select()
